# Apply the "places added" edit to the Terrains worksheet.
#
# Summary of the change:
#  - A new "previousPath" column (D) is introduced.
#  - The old "Forest" terrain row is renamed to "Deep forest" (same
#    enemiesDay/enemiesNight values) and now links back via
#    previousPath to a brand new "Light forest" terrain.
#  - "Cave in the forest" now links back (previousPath) to "Deep forest".
#  - "Cave in the mountains" now links back (previousPath) to "Mountain".
#  - A brand new "Bridge" terrain row is added, which links back
#    (previousPath) to "Riverside".
#  - A final row containing only "Light forest" in column A is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---------------------------------------------------
$ws.Range("D1").Value = "previousPath"

# --- Row 2: Forest -> Deep forest, add previousPath ----------------------
$ws.Range("A2").Value = "Deep forest"
$ws.Range("D2").Value = "Light forest"

# --- Row 5: Cave in the forest, add previousPath --------------------------
$ws.Range("D5").Value = "Deep forest"

# --- Row 9: Cave in the mountains, add previousPath ------------------------
$ws.Range("D9").Value = "Mountain"

# --- Row 10 (new): Bridge --------------------------------------------------
$ws.Range("A10").Value = "Bridge"
$ws.Range("B10").Value = "Mermaid,Troll"
$ws.Range("C10").Value = "Mermaid,Fairy,Troll"
$ws.Range("D10").Value = "Riverside"

# --- Row 11 (new): Light forest (name only) ---------------------------------
$ws.Range("A11").Value = "Light forest"

# --- Column D width / formatting -------------------------------------------
# Mirror Excel's "best fit" auto-sized width for the new column (~18.43
# characters, as produced by double-clicking the column border / AutoFit
# in real Excel once "previousPath" / "Light forest" / "Deep forest" are
# the longest entries).
$ws.Columns.Item(4).AutoFit()
$ws.Columns.Item(4).ColumnWidth = 17.6

# --- Selection matches the authored state -----------------------------------
$ws.Range("D13").Select()
